# verze 4.0.1 opravy v prohlizeci obrazku
$wb = $excel.ActiveWorkbook

# --- Settings sheet ---
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Range("B3").Value = 0
$wsSettings.Range("B5").Value = 1

# --- Settings_recources sheet ---
$wsRes = $wb.Worksheets.Item("Settings_recources")
$wsRes.Range("B3").Value = "C:/Users/jakub.hlavacek.local/Desktop/JHV/test_images/Omron/SnímkySC01_MS500SA/"

# B20/B30 hold numeric-looking text ("40"/"100"); force Text format first so
# Excel keeps them as text instead of auto-converting to numbers.
$wsRes.Range("B20").NumberFormat = "@"
$wsRes.Range("B20").Value = "40"

$wsRes.Range("B30").NumberFormat = "@"
$wsRes.Range("B30").Value = "100"

$wsRes.Range("B31").Value = "ne"
